# Update Excel Mapper sample workbook: allow mapping multiple column names
# (for the same class) across different tabs/files. Sheet1's header row is
# updated to a set of "alternate" column-name variants (suffixed with "1"/"2")
# while Sheet2 keeps the original column names, so the mapper sample can show
# off mapping two different header sets onto the same target class.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1
$ws2 = $wb.Worksheets.Item(2)   # Sheet2

# --- Sheet1: swap in the alternate header names (adds new shared strings) ---
$ws1.Range("A1").Value = "Account Name 1"
$ws1.Range("B1").Value = "Address 2"
$ws1.Range("C1").Value = "My Cool Float 1"
$ws1.Range("D1").Value = "A Decimal Here??? 1"
$ws1.Range("E1").Value = "A simple WhoLeNumber 1"

# --- Sheet1: set portrait page orientation (adds <pageSetup>) ---
$ws1.PageSetup.Orientation = 1

# --- Update selections / active sheet to match the saved view state ---
[void]$ws2.Activate()
[void]$ws2.Range("C39").Select()

[void]$ws1.Activate()
[void]$ws1.Range("E34").Select()
